# Translate the Traditional-Chinese column headers in row 1 to their
# English equivalents (the "座號" and "BMI" headers are left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "height(cm)"
$ws.Range("C1").Value = "weight(kg)"
$ws.Range("E1").Value = "Seated Forward Bend(cm)"
$ws.Range("F1").Value = "Standing long jump(cm)"
$ws.Range("G1").Value = "sit-up(once)"

# Column F now holds the longer "Standing long jump(cm)" header, so widen
# it (originally 16.5 characters wide) to comfortably fit the new text.
$ws.Columns.Item(6).ColumnWidth = 21.86
